$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update final score column (K) for each stock row
$ws.Range("K2").Value = 56.1
$ws.Range("K3").Value = 48.5
$ws.Range("K4").Value = 48.3
$ws.Range("K5").Value = 47.1

# Update MACRO_SCORE column (N) uniformly for all rows
$ws.Range("N2:N5").Value = 53.62998959737769
